$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" values are plain numeric-looking strings
# (e.g. "1.00", "0.0659"). Excel's normal smart-input parsing would
# silently turn those into real numbers (dropping the trailing zero,
# switching to scientific notation, etc.), which would not match the
# original plain-text cell content. For exactly those cells we switch the
# NumberFormat to Text ("@") first so the literal string is kept as-is;
# every other touched cell is left with its original (default) formatting.

$ws.Range("D2").Value = "57.262.36"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "3.095.49"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.78"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.48"
$ws.Range("E6").Value = "  -3.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.096.59"
$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.457"
$ws.Range("E9").Value = "  +2.77%  "

$ws.Range("E10").Value = "  +2.96%  "

$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.399"
$ws.Range("E12").Value = "  +1.72%  "

$ws.Range("D13").Value = "3.627.87"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.40"
$ws.Range("E15").Value = "  -0.90%  "

$ws.Range("E16").Value = "  -1.93%  "

$ws.Range("D17").Value = "57.356.86"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").Value = "3.089.07"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.89"
$ws.Range("E19").Value = "  -3.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.49"
$ws.Range("E20").Value = "  -1.80%  "

$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.45"
$ws.Range("E22").Value = "  +1.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.60"
$ws.Range("E24").Value = "  +1.56%  "

$ws.Range("E25").Value = "  -2.78%  "

$ws.Range("E26").Value = "  -1.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "0.0₃0871"
$ws.Range("E28").Value = "  -6.17%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.87"
$ws.Range("E32").Value = "  -7.59%  "

$ws.Range("E33").Value = "  -0.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("E34").Value = "  +6.15%  "

$ws.Range("E35").Value = "  -3.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.09"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.04"
$ws.Range("E37").Value = "  -2.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.72"
$ws.Range("E38").Value = "  -1.97%  "

$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0659"
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("E41").Value = "  -2.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.04"
$ws.Range("E42").Value = "  +1.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.697"
$ws.Range("E43").Value = "  +1.30%  "

$ws.Range("D44").Value = "2.403.48"
$ws.Range("E44").Value = "  +5.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.69"
$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").Value = "3.135.02"
$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.956"
$ws.Range("E49").Value = "  -3.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.98"
$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.76"
$ws.Range("E51").Value = "  -3.80%  "
